# Scheduled-runner price/profit refresh for the Leve profit sheets
# (ALC, ARM, BSM, CRP, GSM, LTW, WVR). Updates currentAveragePrice /
# LevePrice / LeveProfit columns (H,I,J,K,L,M,N) per-row with refreshed
# market-board figures. CUL is unaffected by this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3239
$ws.Range("J17").Value = 3239
$ws.Range("L17").Value = 9717
$ws.Range("N17").Value = -10053
$ws.Range("H100").Value = 2614.6667
$ws.Range("I100").Value = 2614.6667
$ws.Range("K100").Value = 2614.6667
$ws.Range("M100").Value = -2073.6667
$ws.Range("H112").Value = 1666.1072
$ws.Range("J112").Value = 1677.8889
$ws.Range("L112").Value = 5033.6667
$ws.Range("N112").Value = -7249.6667
$ws.Range("H132").Value = 1852
$ws.Range("I132").Value = 1737.2
$ws.Range("K132").Value = 5211.6
$ws.Range("M132").Value = -2681.6
$ws.Range("H135").Value = 954.5833
$ws.Range("I135").Value = 944.1
$ws.Range("K135").Value = 8496.9
$ws.Range("M135").Value = -5961.9
$ws.Range("H138").Value = 3001.8044
$ws.Range("I138").Value = 1298.3182
$ws.Range("K138").Value = 3894.9546
$ws.Range("M138").Value = 1245.0454
$ws.Range("H140").Value = 106889
$ws.Range("J140").Value = 106889
$ws.Range("L140").Value = 106889
$ws.Range("N140").Value = -117249
$ws.Range("H141").Value = 3233.0527
$ws.Range("I141").Value = 2437
$ws.Range("J141").Value = 9999.5
$ws.Range("K141").Value = 7311
$ws.Range("L141").Value = 29998.5
$ws.Range("M141").Value = -2131
$ws.Range("N141").Value = -40358.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 699.5
$ws.Range("I2").Value = 699.5
$ws.Range("K2").Value = 699.5
$ws.Range("M2").Value = -586.5
$ws.Range("H61").Value = 2032.2222
$ws.Range("I61").Value = 1975.4706
$ws.Range("K61").Value = 1975.4706
$ws.Range("M61").Value = -1763.4706
$ws.Range("H74").Value = 864.2
$ws.Range("I74").Value = 886.5
$ws.Range("J74").Value = 775
$ws.Range("K74").Value = 886.5
$ws.Range("L74").Value = 775
$ws.Range("M74").Value = -12.5
$ws.Range("N74").Value = -2523
$ws.Range("H77").Value = 864.2
$ws.Range("I77").Value = 886.5
$ws.Range("J77").Value = 775
$ws.Range("K77").Value = 4432.5
$ws.Range("L77").Value = 3875
$ws.Range("M77").Value = -64.5
$ws.Range("N77").Value = -12611
$ws.Range("H116").Value = 699.5
$ws.Range("I116").Value = 699.5
$ws.Range("K116").Value = 699.5
$ws.Range("M116").Value = 1594.5
$ws.Range("H132").Value = 2194.457
$ws.Range("I132").Value = 1826.2069
$ws.Range("J132").Value = 3974.3333
$ws.Range("K132").Value = 5478.620699999999
$ws.Range("L132").Value = 11922.9999
$ws.Range("M132").Value = -2948.620699999999
$ws.Range("N132").Value = -16982.9999
$ws.Range("H136").Value = 2032.2222
$ws.Range("I136").Value = 1975.4706
$ws.Range("K136").Value = 5926.4118
$ws.Range("M136").Value = -3376.4118
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 699.5
$ws.Range("I3").Value = 699.5
$ws.Range("K3").Value = 699.5
$ws.Range("M3").Value = -585.5
$ws.Range("H94").Value = 329.125
$ws.Range("I94").Value = 335
$ws.Range("J94").Value = 288
$ws.Range("K94").Value = 335
$ws.Range("L94").Value = 288
$ws.Range("M94").Value = 116
$ws.Range("N94").Value = -1190
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2430.9
$ws.Range("I58").Value = 2314.0625
$ws.Range("J58").Value = 2898.25
$ws.Range("K58").Value = 2314.0625
$ws.Range("L58").Value = 2898.25
$ws.Range("M58").Value = -2111.0625
$ws.Range("N58").Value = -3304.25
$ws.Range("H62").Value = 2597.4
$ws.Range("I62").Value = 2749
$ws.Range("J62").Value = 2496.3333
$ws.Range("K62").Value = 2749
$ws.Range("L62").Value = 2496.3333
$ws.Range("M62").Value = -2125
$ws.Range("N62").Value = -3744.3333
$ws.Range("H65").Value = 2597.4
$ws.Range("I65").Value = 2749
$ws.Range("J65").Value = 2496.3333
$ws.Range("K65").Value = 13745
$ws.Range("L65").Value = 12481.6665
$ws.Range("M65").Value = -10625
$ws.Range("N65").Value = -18721.6665
$ws.Range("H105").Value = 6157
$ws.Range("I105").Value = 5807.636
$ws.Range("K105").Value = 5807.636
$ws.Range("M105").Value = -4060.636
$ws.Range("H132").Value = 3113.6191
$ws.Range("I132").Value = 3102
$ws.Range("J132").Value = 3163
$ws.Range("K132").Value = 9306
$ws.Range("L132").Value = 9489
$ws.Range("M132").Value = -6776
$ws.Range("N132").Value = -14549
$ws.Range("H134").Value = 3798.818
$ws.Range("I134").Value = 3828.75
$ws.Range("J134").Value = 3499.5
$ws.Range("K134").Value = 11486.25
$ws.Range("L134").Value = 10498.5
$ws.Range("M134").Value = -8951.25
$ws.Range("N134").Value = -15568.5
$ws.Range("H136").Value = 2430.9
$ws.Range("I136").Value = 2314.0625
$ws.Range("J136").Value = 2898.25
$ws.Range("K136").Value = 6942.1875
$ws.Range("L136").Value = 8694.75
$ws.Range("M136").Value = -4392.1875
$ws.Range("N136").Value = -13794.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4919.8
$ws.Range("I102").Value = 4919.8
$ws.Range("K102").Value = 4919.8
$ws.Range("M102").Value = -3297.8
$ws.Range("H130").Value = 94966.664
$ws.Range("J130").Value = 94966.664
$ws.Range("L130").Value = 94966.664
$ws.Range("N130").Value = -105006.664
$ws.Range("H132").Value = 3303
$ws.Range("I132").Value = 2917
$ws.Range("J132").Value = 4332.3335
$ws.Range("K132").Value = 8751
$ws.Range("L132").Value = 12997.0005
$ws.Range("M132").Value = -6221
$ws.Range("N132").Value = -18057.0005
$ws.Range("H140").Value = 100700
$ws.Range("J140").Value = 100700
$ws.Range("L140").Value = 100700
$ws.Range("N140").Value = -111060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1956
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 1963.3334
$ws.Range("I100").Value = 1963.3334
$ws.Range("K100").Value = 1963.3334
$ws.Range("M100").Value = -1422.3334
$ws.Range("H132").Value = 2480.818
$ws.Range("I132").Value = 1682.5
$ws.Range("J132").Value = 3438.8
$ws.Range("K132").Value = 5047.5
$ws.Range("L132").Value = 10316.4
$ws.Range("M132").Value = -2517.5
$ws.Range("N132").Value = -15376.4
$ws.Range("H136").Value = 4401.2
$ws.Range("I136").Value = 5001.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 15004.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -12454.5
$ws.Range("N136").Value = -11100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1556.7142
$ws.Range("I122").Value = 1482.8334
$ws.Range("K122").Value = 4448.5002
$ws.Range("M122").Value = -1998.5002
$ws.Range("H126").Value = 2015.2222
$ws.Range("I126").Value = 2167.125
$ws.Range("J126").Value = 800
$ws.Range("K126").Value = 6501.375
$ws.Range("L126").Value = 2400
$ws.Range("M126").Value = -4031.375
$ws.Range("N126").Value = -7340
$ws.Range("H132").Value = 2018.85
$ws.Range("I132").Value = 1552.0667
$ws.Range("K132").Value = 4656.2001
$ws.Range("M132").Value = -2126.2001
$ws.Range("H136").Value = 909.3158
$ws.Range("I136").Value = 909.3158
$ws.Range("K136").Value = 2727.9474
$ws.Range("M136").Value = -177.9474
